$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title area updates (Volume number, date range) ---
$ws.Range("A8").Value2 = "Volume 30   Number  20"
$ws.Range("C9").Value2 = "Report Covering the Week  5/15/2023  Through  5/21/2023"

# --- Crime statistics table updates (rows 14-30) ---

# Cells converting from text placeholder ("0"/"***.*") to real numbers:
# copy NumberFormat from a sibling numeric cell of the desired style first
$fmt15 = $ws.Range("G14").NumberFormat   # style 15 -> #,##0
$fmt16 = $ws.Range("H14").NumberFormat   # style 16 -> #,##0.0 (percent-like)

$ws.Range("D14").NumberFormat = $fmt15
$ws.Range("D14").Value2 = 2
$ws.Range("E14").NumberFormat = $fmt16
$ws.Range("E14").Value2 = -100
$ws.Range("G14").Value2 = 2
$ws.Range("J14").Value2 = 3
$ws.Range("K14").Value2 = -33.333333333333
$ws.Range("F15").Value2 = 4
$ws.Range("C16").Value2 = 8
$ws.Range("D16").Value2 = 2
$ws.Range("E16").Value2 = 300
$ws.Range("F16").Value2 = 27
$ws.Range("G16").Value2 = 18
$ws.Range("H16").Value2 = 50
$ws.Range("I16").Value2 = 116
$ws.Range("J16").Value2 = 93
$ws.Range("K16").Value2 = 24.731182795698
$ws.Range("L16").Value2 = 70.588235294117
$ws.Range("M16").Value2 = 23.404255319148
$ws.Range("N16").Value2 = -62.700964630225
$ws.Range("C17").Value2 = 9
$ws.Range("D17").Value2 = 14
$ws.Range("E17").Value2 = -35.714285714285
$ws.Range("F17").Value2 = 28
$ws.Range("G17").Value2 = 32
$ws.Range("H17").Value2 = -12.5
$ws.Range("I17").Value2 = 143
$ws.Range("J17").Value2 = 152
$ws.Range("K17").Value2 = -5.921052631578
$ws.Range("L17").Value2 = 34.905660377358
$ws.Range("M17").Value2 = -7.142857142857
$ws.Range("N17").Value2 = -11.728395061728
$ws.Range("C18").Value2 = 3
$ws.Range("D18").Value2 = 7
$ws.Range("E18").Value2 = -57.142857142857
$ws.Range("F18").Value2 = 16
$ws.Range("G18").Value2 = 18
$ws.Range("H18").Value2 = -11.111111111111
$ws.Range("I18").Value2 = 75
$ws.Range("J18").Value2 = 63
$ws.Range("K18").Value2 = 19.047619047619
$ws.Range("L18").Value2 = 92.307692307692
$ws.Range("M18").Value2 = -13.793103448275
$ws.Range("N18").Value2 = -78.260869565217
$ws.Range("C19").Value2 = 7
$ws.Range("D19").Value2 = 7
$ws.Range("E19").Value2 = 0
$ws.Range("F19").Value2 = 25
$ws.Range("G19").Value2 = 23
$ws.Range("H19").Value2 = 8.695652173913
$ws.Range("I19").Value2 = 151
$ws.Range("J19").Value2 = 148
$ws.Range("K19").Value2 = 2.027027027027
$ws.Range("L19").Value2 = 37.272727272727
$ws.Range("M19").Value2 = 55.670103092783
$ws.Range("N19").Value2 = 2.721088435374
$ws.Range("C20").Value2 = 6
$ws.Range("E20").Value2 = 0
$ws.Range("F20").Value2 = 22
$ws.Range("G20").Value2 = 18
$ws.Range("H20").Value2 = 22.222222222222
$ws.Range("I20").Value2 = 134
$ws.Range("J20").Value2 = 95
$ws.Range("K20").Value2 = 41.052631578947
$ws.Range("L20").Value2 = 139.285714285714
$ws.Range("M20").Value2 = 362.068965517241
$ws.Range("N20").Value2 = -9.459459459459
$ws.Range("C21").Value2 = 33
$ws.Range("D21").Value2 = 38
$ws.Range("E21").Value2 = -13.157894736842
$ws.Range("F21").Value2 = 122
$ws.Range("G21").Value2 = 111
$ws.Range("H21").Value2 = 9.909909909909
$ws.Range("I21").Value2 = 635
$ws.Range("J21").Value2 = 566
$ws.Range("K21").Value2 = 12.190812720848
$ws.Range("L21").Value2 = 59.949622166246
$ws.Range("M21").Value2 = 35.974304068522
$ws.Range("N21").Value2 = -44.102112676056
$ws.Range("C22").NumberFormat = $fmt15
$ws.Range("C22").Value2 = 3
$ws.Range("D22").NumberFormat = $fmt15
$ws.Range("D22").Value2 = 2
$ws.Range("E22").NumberFormat = $fmt16
$ws.Range("E22").Value2 = 50
$ws.Range("F22").Value2 = 3
$ws.Range("G22").Value2 = 5
$ws.Range("H22").Value2 = -40
$ws.Range("I22").Value2 = 9
$ws.Range("J22").Value2 = 12
$ws.Range("K22").Value2 = -25
$ws.Range("L22").Value2 = 28.571428571428
$ws.Range("M22").Value2 = 12.5
$ws.Range("D23").NumberFormat = $fmt15
$ws.Range("D23").Value2 = 1
$ws.Range("E23").NumberFormat = $fmt16
$ws.Range("E23").Value2 = -100
$ws.Range("G23").NumberFormat = $fmt15
$ws.Range("G23").Value2 = 1
$ws.Range("H23").NumberFormat = $fmt16
$ws.Range("H23").Value2 = -100
$ws.Range("J23").Value2 = 7
$ws.Range("K23").Value2 = 42.857142857142
$ws.Range("C24").Value2 = 14
$ws.Range("E24").Value2 = -6.666666666666
$ws.Range("F24").Value2 = 56
$ws.Range("G24").Value2 = 58
$ws.Range("H24").Value2 = -3.448275862068
$ws.Range("I24").Value2 = 286
$ws.Range("J24").Value2 = 322
$ws.Range("K24").Value2 = -11.180124223602
$ws.Range("L24").Value2 = 60.67415730337
$ws.Range("M24").Value2 = 35.545023696682
$ws.Range("C25").Value2 = 11
$ws.Range("D25").Value2 = 13
$ws.Range("E25").Value2 = -15.384615384615
$ws.Range("F25").Value2 = 36
$ws.Range("G25").Value2 = 52
$ws.Range("H25").Value2 = -30.76923076923
$ws.Range("I25").Value2 = 181
$ws.Range("J25").Value2 = 197
$ws.Range("K25").Value2 = -8.121827411167
$ws.Range("L25").Value2 = 24.827586206896
$ws.Range("M25").Value2 = -18.468468468468
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value2 = "0"
$ws.Range("C26").NumberFormat = "General"
$ws.Range("E26").Value2 = -100
$ws.Range("F26").Value2 = 8
$ws.Range("G26").Value2 = 2
$ws.Range("H26").Value2 = 300
$ws.Range("J26").Value2 = 23
$ws.Range("K26").Value2 = -4.347826086956
$ws.Range("C27").Value2 = 1
$ws.Range("D27").Value2 = 5
$ws.Range("E27").Value2 = -80
$ws.Range("F27").Value2 = 9
$ws.Range("G27").Value2 = 11
$ws.Range("H27").Value2 = -18.181818181818
$ws.Range("I27").Value2 = 41
$ws.Range("J27").Value2 = 51
$ws.Range("K27").Value2 = -19.607843137254
$ws.Range("L27").Value2 = -10.869565217391
$ws.Range("C28").NumberFormat = $fmt15
$ws.Range("C28").Value2 = 3
$ws.Range("D28").NumberFormat = $fmt15
$ws.Range("D28").Value2 = 3
$ws.Range("E28").NumberFormat = $fmt16
$ws.Range("E28").Value2 = 0
$ws.Range("F28").NumberFormat = $fmt15
$ws.Range("F28").Value2 = 3
$ws.Range("G28").Value2 = 3
$ws.Range("H28").Value2 = 0
$ws.Range("I28").Value2 = 9
$ws.Range("J28").Value2 = 9
$ws.Range("L28").Value2 = -10
$ws.Range("M28").Value2 = 28.571428571428
$ws.Range("N28").Value2 = -68.965517241379
$ws.Range("C29").NumberFormat = $fmt15
$ws.Range("C29").Value2 = 3
$ws.Range("D29").NumberFormat = $fmt15
$ws.Range("D29").Value2 = 2
$ws.Range("E29").NumberFormat = $fmt16
$ws.Range("E29").Value2 = 50
$ws.Range("F29").NumberFormat = $fmt15
$ws.Range("F29").Value2 = 3
$ws.Range("G29").Value2 = 2
$ws.Range("H29").Value2 = 50
$ws.Range("I29").Value2 = 8
$ws.Range("J29").Value2 = 7
$ws.Range("K29").Value2 = 14.285714285714
$ws.Range("L29").Value2 = 0
$ws.Range("M29").Value2 = 14.285714285714
$ws.Range("N29").Value2 = -71.428571428571
$ws.Range("D30").NumberFormat = $fmt15
$ws.Range("D30").Value2 = 1
$ws.Range("E30").NumberFormat = $fmt16
$ws.Range("E30").Value2 = -100
$ws.Range("G30").NumberFormat = $fmt15
$ws.Range("G30").Value2 = 1
$ws.Range("H30").NumberFormat = $fmt16
$ws.Range("H30").Value2 = -100
$ws.Range("J30").NumberFormat = $fmt15
$ws.Range("J30").Value2 = 1
$ws.Range("K30").NumberFormat = $fmt16
$ws.Range("K30").Value2 = 0
